$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '68.199.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.51%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.639.49'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.80'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.55%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.542'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("E9").Value = '  -2.17%  '
$ws.Range("E10").Value = '  -1.35%  '
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("E12").Value = '  -1.86%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '27.79'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.83%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.120.78'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.26%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000186'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '68.059.63'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.641.50'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.40%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.38%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '359.77'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.31'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.79%  '
$ws.Range("E21").Value = '  -0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.05'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.28%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '75.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.59%  '
$ws.Range("E25").Value = '  +0.03%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.77'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.51%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.773.44'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.26%  '
$ws.Range("E28").Value = '  -4.59%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '561.84'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.55%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.94'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.81%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.37'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.06%  '
$ws.Range("E33").Value = '  -1.30%  '
$ws.Range("E34").Value = '  -0.04%  '
$ws.Range("E35").Value = '  -3.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -4.15%  '
$ws.Range("E37").Value = '  +0.82%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '158.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.38%  '
$ws.Range("E39").Value = '  -1.88%  '
$ws.Range("E40").Value = '  -2.91%  '
$ws.Range("E41").Value = '  -3.22%  '
$ws.Range("E42").Value = '  +0.35%  '
$ws.Range("E43").Value = '  -4.06%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0₆0321'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -6.01%  '
$ws.Range("E45").Value = '  +0.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '156.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("E47").Value = '  -1.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '21.90'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.96%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.67'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.26%  '
$ws.Range("E50").Value = '  -2.34%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.610'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.25%  '
